# Generate Report for Handoff
# Adds a new handoff record (d9cec067-4e26-449a-b611-f5d7cc4ba428.md) as row 3
# to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$commitSha = "527d1b68398172a2aba8ed92e9803fb038c4c4c2"
$repoBase  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/"
$hoLinkColor = 15570276   # decimal form of RGB 6495ED (matches existing HyperLink font color)

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $hoLinkColor
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
}

function Style-AsDate($rng) {
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Overview sheet: new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "d9cec067-4e26-449a-b611-f5d7cc4ba428.md"
$wsOverview.Range("B3").Value = "e2e\d9cec067-4e26-449a-b611-f5d7cc4ba428.md"
Style-AsHyperlink($wsOverview.Range("B3"))
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 08:42:52"
Style-AsDate($wsOverview.Range("G3"))

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($repoBase + "d9cec067-4e26-449a-b611-f5d7cc4ba428.md"), "", "", "e2e\d9cec067-4e26-449a-b611-f5d7cc4ba428.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "d9cec067-4e26-449a-b611-f5d7cc4ba428.md"
Style-AsHyperlink($wsZhCn.Range("A3"))
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "d9cec067-4e26-449a-b611-f5d7cc4ba428.e163f65a3703e7465144df5113c0811c8318310c.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 08:42:47"
Style-AsDate($wsZhCn.Range("H3"))
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
Style-AsDate($wsZhCn.Range("K3"))
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($repoBase + "d9cec067-4e26-449a-b611-f5d7cc4ba428.md"), "", "", "d9cec067-4e26-449a-b611-f5d7cc4ba428.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "d9cec067-4e26-449a-b611-f5d7cc4ba428.md"
Style-AsHyperlink($wsDeDe.Range("A3"))
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "d9cec067-4e26-449a-b611-f5d7cc4ba428.e163f65a3703e7465144df5113c0811c8318310c.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 08:42:52"
Style-AsDate($wsDeDe.Range("H3"))
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
Style-AsDate($wsDeDe.Range("K3"))
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($repoBase + "d9cec067-4e26-449a-b611-f5d7cc4ba428.md"), "", "", "d9cec067-4e26-449a-b611-f5d7cc4ba428.md") | Out-Null
